# Sync attendance_reports: reorder the "Recorded By" (column G) author list
# so that entries reverse order (e.g. "x, System" -> "System, x") whenever
# "System" appears alongside other recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1 -and ($parts -contains "System")) {
        $n = $parts.Length
        $reversedParts = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversedParts += $parts[$i]
        }
        $cell.Value = $reversedParts -join ", "
    }
}
